$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Remark"

# --- Row 2: Sugar 5KG ---
$ws.Range("A2").Value = "Sugar 5KG"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 400
$ws.Range("D2").Value = "Available"
$ws.Range("E2").Value = "only have 2kg packets. price is Rs.400"

# --- Row 3: biscuits ---
$ws.Range("A3").Value = "biscuits"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 150
$ws.Range("D3").Value = "Available"
$ws.Range("E3").Value = ""

# --- Row 4: Chicken ---
$ws.Range("A4").Value = "Chicken"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 900
$ws.Range("D4").Value = "Available"
$ws.Range("E4").Value = ""

# --- Remove the old Status/Remark columns (F, G) which are now redundant ---
$ws.Range("F1:G1").EntireColumn.Delete()

# --- Resize the new Status/Remark columns (D, E) ---
$ws.Range("D1").EntireColumn.ColumnWidth = 21.5
$ws.Range("E1").EntireColumn.ColumnWidth = 32.5
